$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.207.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.602.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '2.32'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +20.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '226.36'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '636.37'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.411'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.599.81'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.48'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.206'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000289'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.274.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.981.03'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.75'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.47'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.606.34'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.512'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '505.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.23'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.37%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +23.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '117.77'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +16.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000203'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.74'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.60'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.74'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.90'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.179'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.78'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.75'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.583'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '591.33'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.30'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.78'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.77'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.473'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0472'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.91'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.917'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.45'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.57'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'MantraDAO'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.58'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.20'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.36%  '